$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.946.99'
$ws.Range('D3').Value = '1.671.88'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '214.88'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '20.16'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').Value = '1.907.15'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '1.673.00'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '65.55'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = '26.945.21'
$ws.Range('D18').Value = '234.28'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '8.04'
$ws.Range('E19').Value = '  +3.37%  '
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = '9.16'
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('D26').Value = '7.14'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').Value = '15.97'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('D32').Value = '3.33'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').Value = '1.465.19'
$ws.Range('E33').Value = '  -5.25%  '
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('D35').Value = '1.66'
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').Value = '0.579'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').Value = '0.897'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '0.0171'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('E40').Value = '  +12.57%  '
$ws.Range('D41').Value = '5.78'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.30'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.70'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = '1.811.86'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = '0.779'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '90.65'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').Value = '7.69'
$ws.Range('E51').Value = '  +0.15%  '
